# Fix "Rougelike" -> "Roguelike" typo across the Action Roguelike project
# row (name, banner image path, and the content cell's embedded image
# paths / bullet text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Project name cell (A5)
$ws.Range("A5").Value = "Action Roguelike"

# Banner image path cell (B5)
$ws.Range("B5").Value = "/images/projects/action-roguelike/ar-banner.png"

# Content cell (F5): fix the "Rougelike combat" bullet text and the two
# embedded "action-rougelike" image paths, keep everything else identical.
$f5 = @"
3#https://youtu.be/mfHQZCE7S80
0#This was my university graduation project in my final year.
0# What I worked on:\\n- Roguelike combat. \\n- Random rooms and upgrades system. \\n- Enemy AI, trap system. \\n- Save load system.\\n- Event trigger, simple shop.
0# Full gameplay.
3#https://youtu.be/8IbrqM4H0rg
1#/images/projects/action-roguelike/ar-combat-2.png
2#/images/projects/action-roguelike/ar-combat.gif
"@

$ws.Range("F5").Value = $f5

# Setting a wrapped multi-line cell's value can trigger an auto row-height
# recalculation in some engines; the line count here is unchanged from the
# original, so restore the row's original display height.
$ws.Rows(5).RowHeight = 18

# Leave the selection where the editing session ended up.
[void]$ws.Range("C26").Select()
